# Auto-generated data refresh: updates market-price derived columns (H:N)
# across multiple Leve-profit worksheets. Values are plain numbers (no formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3087.0435
$ws.Range("I64").Value = 3095.3333
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 3095.3333
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2847.3333
$ws.Range("N64").Value = -3496

# Row 67
$ws.Range("H67").Value = 3087.0435
$ws.Range("I67").Value = 3095.3333
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 3095.3333
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2237.3333
$ws.Range("N67").Value = -4716

# Row 107
$ws.Range("H107").Value = 9222.308000000001
$ws.Range("I107").Value = 9222.308000000001
$ws.Range("K107").Value = 9222.308000000001
$ws.Range("M107").Value = -7302.308000000001

# Row 137
$ws.Range("H137").Value = 1223.3334
$ws.Range("I137").Value = 1347.25
$ws.Range("J137").Value = 1124.2
$ws.Range("K137").Value = 4041.75
$ws.Range("L137").Value = 3372.6
$ws.Range("M137").Value = -1491.75
$ws.Range("N137").Value = -8472.6

# Row 138
$ws.Range("H138").Value = 2596.9834
$ws.Range("I138").Value = 3718.6
$ws.Range("J138").Value = 2372.66
$ws.Range("K138").Value = 11155.8
$ws.Range("L138").Value = 7117.98
$ws.Range("M138").Value = -6015.799999999999
$ws.Range("N138").Value = -17397.98

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1314.762
$ws.Range("I2").Value = 1256.0625
$ws.Range("J2").Value = 1502.6
$ws.Range("K2").Value = 1256.0625
$ws.Range("L2").Value = 1502.6
$ws.Range("M2").Value = -1143.0625
$ws.Range("N2").Value = -1728.6

# Row 32
$ws.Range("H32").Value = 560421.75
$ws.Range("I32").Value = 615463.8
$ws.Range("K32").Value = 615463.8
$ws.Range("M32").Value = -615176.8

# Row 97
$ws.Range("H97").Value = 1146.7354
$ws.Range("I97").Value = 1131.12
$ws.Range("K97").Value = 1131.12
$ws.Range("M97").Value = -635.1199999999999

# Row 116
$ws.Range("H116").Value = 1314.762
$ws.Range("I116").Value = 1256.0625
$ws.Range("J116").Value = 1502.6
$ws.Range("K116").Value = 1256.0625
$ws.Range("L116").Value = 1502.6
$ws.Range("M116").Value = 1037.9375
$ws.Range("N116").Value = -6090.6

# Row 132
$ws.Range("H132").Value = 3710.9714
$ws.Range("I132").Value = 4223.2383
$ws.Range("J132").Value = 2942.5715
$ws.Range("K132").Value = 12669.7149
$ws.Range("L132").Value = 8827.7145
$ws.Range("M132").Value = -10139.7149
$ws.Range("N132").Value = -13887.7145

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1314.762
$ws.Range("I3").Value = 1256.0625
$ws.Range("J3").Value = 1502.6
$ws.Range("K3").Value = 1256.0625
$ws.Range("L3").Value = 1502.6
$ws.Range("M3").Value = -1142.0625
$ws.Range("N3").Value = -1730.6

# Row 35
$ws.Range("H35").Value = 19999
$ws.Range("J35").Value = 19999
$ws.Range("L35").Value = 19999
$ws.Range("N35").Value = -20619

# Row 86
$ws.Range("H86").Value = 1953.3334
$ws.Range("I86").Value = 1953.3334
$ws.Range("K86").Value = 1953.3334
$ws.Range("M86").Value = -830.3334

# Row 89
$ws.Range("H89").Value = 1953.3334
$ws.Range("I89").Value = 1953.3334
$ws.Range("K89").Value = 9766.666999999999
$ws.Range("M89").Value = -4150.666999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6128.1377
$ws.Range("I31").Value = 1606.6
$ws.Range("J31").Value = 10972.643
$ws.Range("K31").Value = 1606.6
$ws.Range("L31").Value = 10972.643
$ws.Range("M31").Value = -1311.6
$ws.Range("N31").Value = -11562.643

# Row 34
$ws.Range("H34").Value = 6128.1377
$ws.Range("I34").Value = 1606.6
$ws.Range("J34").Value = 10972.643
$ws.Range("K34").Value = 1606.6
$ws.Range("L34").Value = 10972.643
$ws.Range("M34").Value = -1404.6
$ws.Range("N34").Value = -11376.643

# Row 58
$ws.Range("H58").Value = 2560.2917
$ws.Range("I58").Value = 2379.2942
$ws.Range("J58").Value = 2999.8572
$ws.Range("K58").Value = 2379.2942
$ws.Range("L58").Value = 2999.8572
$ws.Range("M58").Value = -2176.2942
$ws.Range("N58").Value = -3405.8572

# Row 62
$ws.Range("H62").Value = 4075.4546
$ws.Range("I62").Value = 4075.4546
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4075.4546
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3451.4546
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 4075.4546
$ws.Range("I65").Value = 4075.4546
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20377.273
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -17257.273
$ws.Range("N65").ClearContents()

# Row 132
$ws.Range("H132").Value = 20837272
$ws.Range("I132").Value = 4333.3335
$ws.Range("J132").Value = 33337034
$ws.Range("K132").Value = 13000.0005
$ws.Range("L132").Value = 100011102
$ws.Range("M132").Value = -10470.0005
$ws.Range("N132").Value = -100016162

# Row 133
$ws.Range("H133").Value = 47990
$ws.Range("J133").Value = 47990
$ws.Range("L133").Value = 47990
$ws.Range("N133").Value = -53050

# Row 136
$ws.Range("H136").Value = 2560.2917
$ws.Range("I136").Value = 2379.2942
$ws.Range("J136").Value = 2999.8572
$ws.Range("K136").Value = 7137.882599999999
$ws.Range("L136").Value = 8999.571599999999
$ws.Range("M136").Value = -4587.882599999999
$ws.Range("N136").Value = -14099.5716

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 4000
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = -5831
$ws.Range("N17").Value = -12338

# Row 19
$ws.Range("H19").Value = 5500
$ws.Range("I19").Value = 6500
$ws.Range("K19").Value = 19500
$ws.Range("M19").Value = -19326

# Row 22
$ws.Range("H22").Value = 1700
$ws.Range("I22").Value = 920
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 2760
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = -2591
$ws.Range("N22").Value = -9338

# Row 27
$ws.Range("H27").Value = 1700
$ws.Range("I27").Value = 920
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 2760
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = -2658
$ws.Range("N27").Value = -9204

# Row 80
$ws.Range("H80").Value = 4916.6665
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 5400
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 16200
$ws.Range("M80").Value = -6564
$ws.Range("N80").Value = -18072

# Row 83
$ws.Range("H83").Value = 4916.6665
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 5400
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 48600
$ws.Range("M83").Value = -17820
$ws.Range("N83").Value = -57960

# Row 137
$ws.Range("H137").Value = 6849.5557
$ws.Range("J137").Value = 4670
$ws.Range("L137").Value = 14010
$ws.Range("N137").Value = -24210

$ws = $wb.Worksheets.Item("GSM")
# Row 69
$ws.Range("H69").Value = 39285.715
$ws.Range("J69").Value = 39285.715
$ws.Range("L69").Value = 39285.715
$ws.Range("N69").Value = -40783.715

# Row 70
$ws.Range("H70").Value = 5697.1387
$ws.Range("I70").Value = 5656.32
$ws.Range("J70").Value = 5789.909
$ws.Range("K70").Value = 5656.32
$ws.Range("L70").Value = 5789.909
$ws.Range("M70").Value = -5386.32
$ws.Range("N70").Value = -6329.909

# Row 72
$ws.Range("H72").Value = 39285.715
$ws.Range("J72").Value = 39285.715
$ws.Range("L72").Value = 117857.145
$ws.Range("N72").Value = -125345.145

# Row 73
$ws.Range("H73").Value = 5697.1387
$ws.Range("I73").Value = 5656.32
$ws.Range("J73").Value = 5789.909
$ws.Range("K73").Value = 5656.32
$ws.Range("L73").Value = 5789.909
$ws.Range("M73").Value = -4720.32
$ws.Range("N73").Value = -7661.909

# Row 132
$ws.Range("H132").Value = 2155.8286
$ws.Range("I132").Value = 1873.4642
$ws.Range("J132").Value = 3285.2856
$ws.Range("K132").Value = 5620.392599999999
$ws.Range("L132").Value = 9855.856800000001
$ws.Range("M132").Value = -3090.392599999999
$ws.Range("N132").Value = -14915.8568

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 6028.857
$ws.Range("I122").Value = 3976
$ws.Range("K122").Value = 11928
$ws.Range("M122").Value = -9478

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 18334.5
$ws.Range("J15").Value = 18334.5
$ws.Range("L15").Value = 18334.5
$ws.Range("N15").Value = -18910.5

# Row 82
$ws.Range("H82").Value = 69140.5
$ws.Range("J82").Value = 69140.5
$ws.Range("L82").Value = 69140.5
$ws.Range("N82").Value = -69906.5

# Row 85
$ws.Range("H85").Value = 69140.5
$ws.Range("J85").Value = 69140.5
$ws.Range("L85").Value = 69140.5
$ws.Range("N85").Value = -71792.5

# Row 96
$ws.Range("H96").Value = 6695
$ws.Range("I96").Value = 4090
$ws.Range("J96").Value = 9300
$ws.Range("K96").Value = 4090
$ws.Range("L96").Value = 9300
$ws.Range("M96").Value = -2717
$ws.Range("N96").Value = -12046

# Row 122
$ws.Range("H122").Value = 2100.889
$ws.Range("I122").Value = 1421.6
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 4264.799999999999
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -1814.799999999999
$ws.Range("N122").Value = -13750

# Row 132
$ws.Range("H132").Value = 8335713.5
$ws.Range("I132").Value = 2691.6428
$ws.Range("J132").Value = 13891062
$ws.Range("K132").Value = 8074.928400000001
$ws.Range("L132").Value = 41673186
$ws.Range("M132").Value = -5544.928400000001
$ws.Range("N132").Value = -41678246

# Row 136
$ws.Range("H136").Value = 3089.4614
$ws.Range("I136").Value = 2912.25
$ws.Range("J136").Value = 3373
$ws.Range("K136").Value = 8736.75
$ws.Range("L136").Value = 10119
$ws.Range("M136").Value = -6186.75
$ws.Range("N136").Value = -15219
